$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1803.1333
$ws.Range("I106").Value = 1520.5834
$ws.Range("K106").Value = 1520.5834
$ws.Range("M106").Value = -889.5834
$ws.Range("H132").Value = 3381548
$ws.Range("I132").Value = 3475462
$ws.Range("J132").Value = 644
$ws.Range("K132").Value = 10426386
$ws.Range("L132").Value = 1932
$ws.Range("M132").Value = -10423856
$ws.Range("N132").Value = -6992
$ws.Range("H137").Value = 1742.2285
$ws.Range("I137").Value = 1364.4615
$ws.Range("J137").Value = 2833.5557
$ws.Range("K137").Value = 4093.3845
$ws.Range("L137").Value = 8500.667099999999
$ws.Range("M137").Value = -1543.3845
$ws.Range("N137").Value = -13600.6671
$ws.Range("H138").Value = 1474.8955
$ws.Range("I138").Value = 991.37256
$ws.Range("J138").Value = 3016.125
$ws.Range("K138").Value = 2974.11768
$ws.Range("L138").Value = 9048.375
$ws.Range("M138").Value = 2165.88232
$ws.Range("N138").Value = -19328.375
$ws.Range("H141").Value = 2090.7896
$ws.Range("I141").Value = 1970.9722
$ws.Range("J141").Value = 4247.5
$ws.Range("K141").Value = 5912.9166
$ws.Range("L141").Value = 12742.5
$ws.Range("M141").Value = -732.9165999999996
$ws.Range("N141").Value = -23102.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -188
$ws.Range("N5").Value = ""
$ws.Range("H32").Value = 30758.95
$ws.Range("I32").Value = 6317.7925
$ws.Range("K32").Value = 6317.7925
$ws.Range("M32").Value = -6030.7925
$ws.Range("H45").Value = 251323
$ws.Range("I45").Value = 251323
$ws.Range("K45").Value = 251323
$ws.Range("M45").Value = -250946
$ws.Range("H61").Value = 2405.6572
$ws.Range("I61").Value = 1449.9166
$ws.Range("J61").Value = 2904.3044
$ws.Range("K61").Value = 1449.9166
$ws.Range("L61").Value = 2904.3044
$ws.Range("M61").Value = -1237.9166
$ws.Range("N61").Value = -3328.3044
$ws.Range("H74").Value = 684
$ws.Range("I74").Value = 508
$ws.Range("J74").Value = 816
$ws.Range("K74").Value = 508
$ws.Range("L74").Value = 816
$ws.Range("M74").Value = 366
$ws.Range("N74").Value = -2564
$ws.Range("H77").Value = 684
$ws.Range("I77").Value = 508
$ws.Range("J77").Value = 816
$ws.Range("K77").Value = 2540
$ws.Range("L77").Value = 4080
$ws.Range("M77").Value = 1828
$ws.Range("N77").Value = -12816
$ws.Range("H132").Value = 2608.4565
$ws.Range("I132").Value = 2854.861
$ws.Range("J132").Value = 1721.4
$ws.Range("K132").Value = 8564.582999999999
$ws.Range("L132").Value = 5164.200000000001
$ws.Range("M132").Value = -6034.582999999999
$ws.Range("N132").Value = -10224.2
$ws.Range("H136").Value = 2405.6572
$ws.Range("I136").Value = 1449.9166
$ws.Range("J136").Value = 2904.3044
$ws.Range("K136").Value = 4349.7498
$ws.Range("L136").Value = 8712.913199999999
$ws.Range("M136").Value = -1799.7498
$ws.Range("N136").Value = -13812.9132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -185
$ws.Range("N4").Value = ""
$ws.Range("H134").Value = 1828.6744
$ws.Range("I134").Value = 1828.6744
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5486.023200000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2951.023200000001
$ws.Range("N134").Value = ""
$ws.Range("H141").Value = 36944.25
$ws.Range("J141").Value = 45925.668
$ws.Range("L141").Value = 45925.668
$ws.Range("N141").Value = -56285.668

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30638.24
$ws.Range("I31").Value = 1997.8422
$ws.Range("J31").Value = 48192.03
$ws.Range("K31").Value = 1997.8422
$ws.Range("L31").Value = 48192.03
$ws.Range("M31").Value = -1702.8422
$ws.Range("N31").Value = -48782.03
$ws.Range("H34").Value = 30638.24
$ws.Range("I34").Value = 1997.8422
$ws.Range("J34").Value = 48192.03
$ws.Range("K34").Value = 1997.8422
$ws.Range("L34").Value = 48192.03
$ws.Range("M34").Value = -1795.8422
$ws.Range("N34").Value = -48596.03
$ws.Range("H58").Value = 943.3090999999999
$ws.Range("I58").Value = 873.8077
$ws.Range("J58").Value = 2148
$ws.Range("K58").Value = 873.8077
$ws.Range("L58").Value = 2148
$ws.Range("M58").Value = -670.8077
$ws.Range("N58").Value = -2554
$ws.Range("H107").Value = 758.94116
$ws.Range("I107").Value = 898
$ws.Range("J107").Value = 560.2857
$ws.Range("K107").Value = 898
$ws.Range("L107").Value = 560.2857
$ws.Range("M107").Value = 1022
$ws.Range("N107").Value = -4400.2857
$ws.Range("H132").Value = 23079236
$ws.Range("I132").Value = 20002138
$ws.Range("J132").Value = 33336228
$ws.Range("K132").Value = 60006414
$ws.Range("L132").Value = 100008684
$ws.Range("M132").Value = -60003884
$ws.Range("N132").Value = -100013744
$ws.Range("H134").Value = 900.2692
$ws.Range("I134").Value = 751.70734
$ws.Range("K134").Value = 2255.12202
$ws.Range("M134").Value = 279.8779799999998
$ws.Range("H136").Value = 943.3090999999999
$ws.Range("I136").Value = 873.8077
$ws.Range("J136").Value = 2148
$ws.Range("K136").Value = 2621.4231
$ws.Range("L136").Value = 6444
$ws.Range("M136").Value = -71.42309999999998
$ws.Range("N136").Value = -11544

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 7664.143
$ws.Range("J22").Value = 8841.5
$ws.Range("L22").Value = 26524.5
$ws.Range("N22").Value = -26862.5
$ws.Range("H27").Value = 7664.143
$ws.Range("J27").Value = 8841.5
$ws.Range("L27").Value = 26524.5
$ws.Range("N27").Value = -26728.5
$ws.Range("H121").Value = 7725.1333
$ws.Range("I121").Value = 7024.7144
$ws.Range("J121").Value = 8338
$ws.Range("K121").Value = 21074.1432
$ws.Range("L121").Value = 25014
$ws.Range("M121").Value = -19764.1432
$ws.Range("N121").Value = -27634

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 47762.4
$ws.Range("J19").Value = 47762.4
$ws.Range("L19").Value = 47762.4
$ws.Range("N19").Value = -48338.4
$ws.Range("H113").Value = 1457.6154
$ws.Range("I113").Value = 1233.25
$ws.Range("J113").Value = 1816.6
$ws.Range("K113").Value = 1233.25
$ws.Range("L113").Value = 1816.6
$ws.Range("M113").Value = 936.75
$ws.Range("N113").Value = -6156.6
$ws.Range("H132").Value = 2469.5715
$ws.Range("I132").Value = 2547.3513
$ws.Range("J132").Value = 1894
$ws.Range("K132").Value = 7642.053899999999
$ws.Range("L132").Value = 5682
$ws.Range("M132").Value = -5112.053899999999
$ws.Range("N132").Value = -10742

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 863.2549
$ws.Range("I136").Value = 666.775
$ws.Range("J136").Value = 1577.7273
$ws.Range("K136").Value = 2000.325
$ws.Range("L136").Value = 4733.1819
$ws.Range("M136").Value = 549.6750000000002
$ws.Range("N136").Value = -9833.1819

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 140400
$ws.Range("I26").Value = 2000
$ws.Range("J26").Value = 186533.33
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 186533.33
$ws.Range("M26").Value = -1707
$ws.Range("N26").Value = -187119.33
$ws.Range("H132").Value = 1264.3334
$ws.Range("I132").Value = 1294.5363
$ws.Range("K132").Value = 3883.6089
$ws.Range("M132").Value = -1353.6089
